$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap match-result/odds data (columns F:V) between row pairs. ---
# The scraper re-ordered these fixture rows; index/date columns (A:E) stay put.
function Swap-Rows($r1, $r2) {
    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

$swapPairs = @(
    @(26, 27),
    @(70, 71),
    @(72, 73),
    @(80, 81),
    @(86, 87),
    @(89, 90),
    @(96, 97)
)

foreach ($pair in $swapPairs) {
    Swap-Rows $pair[0] $pair[1]
}

# --- 2. Append newly scraped fixtures as rows 117-121. ---
# Clone formatting (styles s=1 / s=2 on columns A and E) from the last
# existing data row before filling in the new values.
$ws.Range("A116:V116").Copy()
$ws.Range("A117:V121").PasteSpecial(-4122)

$newRows = @(
    @{ RowNum=117; A=116; E=45242.66666666666; F="Celta Vigo B";     G=2; H="Unionistas";        I=1; J=1.72; K="09/11/2023 09:13"; L=1.94; M="12/11/2023 15:58"; N=3.37; O="09/11/2023 09:13"; P=3.43; Q="12/11/2023 15:58"; R=4.44; S="09/11/2023 09:13"; T=3.98; U="12/11/2023 15:58"; V="https://www.betexplorer.com/football/spain/primera-rfef-group-1/celta-vigo-unionistas-de-salamanca/ILQo7oL1/" },
    @{ RowNum=118; A=117; E=45242.66666666666; F="Fuenlabrada";      G=4; H="Osasuna B";         I=1; J=2.15; K="09/11/2023 09:13"; L=2.06; M="12/11/2023 15:51"; N=3;    O="09/11/2023 09:13"; P=3.23; Q="12/11/2023 15:51"; R=3.3;  S="09/11/2023 09:13"; T=3.8;  U="12/11/2023 15:51"; V="https://www.betexplorer.com/football/spain/primera-rfef-group-1/cf-fuenlabrada-osasuna/WOUk65z8/" },
    @{ RowNum=119; A=118; E=45242.70833333334; F="Lugo";             G=0; H="Leonesa";           I=3; J=2.14; K="09/11/2023 09:13"; L=2.7;  M="12/11/2023 16:42"; N=2.93; O="09/11/2023 09:13"; P=2.86; Q="12/11/2023 16:42"; R=3.36; S="09/11/2023 09:13"; T=2.97; U="12/11/2023 16:42"; V="https://www.betexplorer.com/football/spain/primera-rfef-group-1/lugo-leonesa/IJZ9B0Cm/" },
    @{ RowNum=120; A=119; E=45242.75;          F="Ponferradina";     G=1; H="Rayo Majadahonda";  I=0; J=1.6;  K="09/11/2023 09:13"; L=1.63; M="12/11/2023 13:47"; N=3.55; O="09/11/2023 09:13"; P=3.81; Q="12/11/2023 16:02"; R=5.04; S="09/11/2023 09:13"; T=5.44; U="12/11/2023 13:47"; V="https://www.betexplorer.com/football/spain/primera-rfef-group-1/ponferradina-rayo-majadahonda/vTz6Ctcs/" },
    @{ RowNum=121; A=120; E=45242.83333333334; F="Gimnastic";        G=0; H="Sestao";            I=0; J=1.64; K="09/11/2023 09:13"; L=1.7;  M="12/11/2023 19:50"; N=3.48; O="09/11/2023 09:13"; P=3.48; Q="12/11/2023 19:50"; R=4.94; S="09/11/2023 09:13"; T=5.46; U="12/11/2023 19:50"; V="https://www.betexplorer.com/football/spain/primera-rfef-group-1/gimnastic-de-tarragona-sestao/d6xw97jk/" }
)

foreach ($row in $newRows) {
    $r = $row.RowNum
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = "spain"
    $ws.Range("C$r").Value = "primera-rfef-group-1"
    $ws.Range("D$r").Value = "2023-2024"
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $ws.Range("N$r").Value = $row.N
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    $ws.Range("Q$r").Value = $row.Q
    $ws.Range("R$r").Value = $row.R
    $ws.Range("S$r").Value = $row.S
    $ws.Range("T$r").Value = $row.T
    $ws.Range("U$r").Value = $row.U
    $ws.Range("V$r").Value = $row.V
}

"done"
